$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.563.42"
$ws.Range("E2").Value = "  -6.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.552.41"
$ws.Range("E3").Value = "  -4.69%  "

$ws.Range("E4").Value = "  -0.56%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "395.00"
$ws.Range("E5").Value = "  -6.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "122.92"
$ws.Range("E6").Value = "  -6.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.541.25"
$ws.Range("E7").Value = "  -4.67%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  -9.67%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.681"
$ws.Range("E10").Value = "  -11.42%  "

$ws.Range("E11").Value = "  -20.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000328"
$ws.Range("E12").Value = "  -20.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.94"
$ws.Range("E13").Value = "  -8.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.109.70"
$ws.Range("E14").Value = "  -4.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.17"
$ws.Range("E15").Value = "  -7.45%  "

$ws.Range("E16").Value = "  -2.86%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.79"
$ws.Range("E17").Value = "  +8.24%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.542.06"
$ws.Range("E18").Value = "  -4.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.66"
$ws.Range("E19").Value = "  -9.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "63.492.70"
$ws.Range("E20").Value = "  -6.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.02"
$ws.Range("E21").Value = "  -11.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "393.85"
$ws.Range("E22").Value = "  -13.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.05"
$ws.Range("E23").Value = "  +1.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.91"
$ws.Range("E24").Value = "  -8.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.90"
$ws.Range("E25").Value = "  -5.88%  "

$ws.Range("E26").Value = "  +9.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.97"
$ws.Range("E27").Value = "  -13.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.99"
$ws.Range("E28").Value = "  -8.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.72"
$ws.Range("E29").Value = "  -14.05%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.96"
$ws.Range("E30").Value = "  -3.19%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.60"
$ws.Range("E31").Value = "  -7.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.113"
$ws.Range("E32").Value = "  -5.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.74"
$ws.Range("E33").Value = "  -5.87%  "

$ws.Range("E34").Value = "  -6.29%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "36.60"
$ws.Range("E36").Value = "  -9.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.65"
$ws.Range("E37").Value = "  -3.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0437"
$ws.Range("E38").Value = "  -10.73%  "

$ws.Range("E39").Value = "  -0.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0653"
$ws.Range("E40").Value = "  -13.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.68"
$ws.Range("E41").Value = "  +3.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.131"
$ws.Range("E42").Value = "  -12.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.08"
$ws.Range("E43").Value = "  +16.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "141.13"
$ws.Range("E44").Value = "  -5.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.69"
$ws.Range("E45").Value = "  +10.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.96"
$ws.Range("E46").Value = "  -6.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.07"
$ws.Range("E47").Value = "  -10.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.08"
$ws.Range("E48").Value = "  -4.96%  "

$ws.Range("E49").Value = "  -8.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.64"
$ws.Range("E50").Value = "  -9.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.275"
$ws.Range("E51").Value = "  -9.10%  "

